$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = "Fase 4"

$ws.Range("B28").Value = "Broadcast"
$ws.Range("C28").Value = 1

$ws.Range("B29").Value = "Valgrind"
$ws.Range("C29").Value = 10

$ws.Range("B30").Value = "Debuggar"
$ws.Range("C30").Value = 20

$ws.Range("B31").Value = "Memòria"
$ws.Range("C31").Value = 12

$ws.Range("C31").Select()
